# Appends the six new paragraphs (two empty yellow-highlighted marks,
# then the Liao 2018 / Perez 2021 / Toussaint 2009 / Gabana 2017 entries)
# to the end of the document body, right before the final section break,
# exactly reproducing the author's new content/formatting.

$d = $word.ActiveDocument

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:lastRenderedPageBreak/><w:t>Liao 2018:</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t>Examining the mechanism between gratefulness and SWB.</w:t></w:r><w:r><w:t xml:space="preserve"> Working </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>off of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ‘broaden and build theory’, the hypothesized mediators are social and cognitive resources. Mediation was tested using latent change score analysis and SEM techniques, two mediators ‘social connectedness’ and ‘presence of meaning in life’ were significant mediators.</w:t></w:r><w:r><w:t xml:space="preserve"> Tested using Shrout and Bolger’s bootstrap method</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Perez 2021</w:t></w:r><w:r><w:t>: Mediation of spirituality on life satisfaction</w:t></w:r><w:r><w:t xml:space="preserve"> in teens/young adults from the Philippines</w:t></w:r><w:r><w:t xml:space="preserve">, the majority of whom are religious. </w:t></w:r><w:r><w:t>Mediation done using the Hayes bootstrapping method. Gratitude significantly affected spirituality and life satisfaction, spirituality directly affected life satisfaction, and there was a mediating effect of spirituality on life satisfaction!</w:t></w:r><w:r><w:t xml:space="preserve"> Partial mediation.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Toussaint 2009:</w:t></w:r><w:r><w:t xml:space="preserve">  Forgiveness and gratitude both associated with SWB, and this is partially mediated by affect and belief, done on 72</w:t></w:r><w:r><w:t xml:space="preserve"> clinical psychological</w:t></w:r><w:r><w:t xml:space="preserve"> outpatients in </w:t></w:r><w:r><w:t>Philadelphia</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Used Baron and Kenny procedure to test mediation.</w:t></w:r><w:r><w:t xml:space="preserve"> Positive affect was a significant mediator, as well as ‘belief in self’/’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>self worth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’. We saw full or near-full mediation.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Gabana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> 2017:</w:t></w:r><w:r><w:t xml:space="preserve"> Study done on D1 and D3 college athletes, relationship between gratitude and </w:t></w:r><w:r><w:t xml:space="preserve">athletic </w:t></w:r><w:r><w:t>burnout, as well as gratitude and sport satisfaction, both were mediated strongly by ‘perceived social support’</w:t></w:r><w:r><w:t xml:space="preserve">. Defined </w:t></w:r><w:bookmarkStart w:id="1" w:name="_Hlk116568768"/><w:r><w:t>as “</w:t></w:r><w:r><w:t>one’s potential</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>access to social support and is a support recipient’s subjective judgment that friends,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>family, team-mates, and coaches would provide assistance if needed</w:t></w:r><w:r><w:t>”</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p>
'@

$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertXML($newParagraphsXml)

Write-Host ("Paragraph count after insert: " + $d.Paragraphs.Count)
